# Auto-generated Excel COM-interop script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("M1").Value = "photoOriginalURL"
$ws.Range("N1").Value = "photoFileName"
$ws.Range("O1").Value = "photoNewName"
$ws.Range("P1").Value = "photoNewURL"
$ws.Range("Q1").Value = "fileNameConversionScript"
$ws.Range("R1").Value = "articleIntro"
$ws.Range("S1").Value = "photoNewFileNameFull"
$ws.Range("T1").Value = "structuredData"
$ws.Range("U1").Value = "photoLocal"
$ws.Range("V1").Value = "missingData"

# ---- Row 2 (Vivero Santa Fe) ----
$ws.Range("N2").Value = 'AF1QipOC7w7iQu2kg7BiEaF1hoYdzwUX11Mfy6wYJIgK=w408-h306-k-no'
$ws.Range("O2").Value = 'Vivero_Santa_Fe'
$ws.Range("P2").Value = 'https://rumbonaturaleza.com/uploads/2023/03/Vivero_Santa_Fe.jpg'
$ws.Range("Q2").Value = 'ren "AF1QipOC7w7iQu2kg7BiEaF1hoYdzwUX11Mfy6wYJIgK=w408-h306-k-no" "Vivero_Santa_Fe.jpg"'
$r2 = @'

                    <p> ¿Estás buscando los mejores Parques Ecoturísticos en Tabasco? ¡Estás en el lugar correcto! Pues en este artículo vamos a presentarte cuáles son los  Parque Ecoturístico que han sido mejor evaluados en este estado. 
 Para esto, realizamos consultas en un montón de fuentes oficiales, redes sociales, rankings e incluso entrevistas para poder determinar cuáles son los  Parque Ecoturístico que mejor calificación han recibido en Tabasco durante los últimos años. 
 Con esta prueba social como respaldo, hoy te daremos los Parque Ecoturístico mejor calificados y te compartiremos su ubicación, medios oficiales de contacto, horarios y cómo llegar hasta ellos, junto con la calificación promedio con la que cuenta cada lugar. 
 Así que prepárate y ¡a disfrutar del ecoturismo en Tabasco!</p>                    
                    
'@
$ws.Range("R2").Value = $r2
$s2 = @'

                    <p> ¿Estás buscando los mejores Parques Ecoturísticos en Tabasco? ¡Estás en el lugar correcto! Pues en este artículo vamos a presentarte cuáles son los  Parque Ecoturístico que han sido mejor evaluados en este estado. 
 Para esto, realizamos consultas en un montón de fuentes oficiales, redes sociales, rankings e incluso entrevistas para poder determinar cuáles son los  Parque Ecoturístico que mejor calificación han recibido en Tabasco durante los últimos años. 
 Con esta prueba social como respaldo, hoy te daremos los Parque Ecoturístico mejor calificados y te compartiremos su ubicación, medios oficiales de contacto, horarios y cómo llegar hasta ellos, junto con la calificación promedio con la que cuenta cada lugar. 
 Así que prepárate y ¡a disfrutar del ecoturismo en Tabasco!</p>                    
                    
'@
$ws.Range("S2").Value = $s2
$t2 = @'

                    <h2><b>Parque Ecoturístico Vivero Santa Fe</b></h2>
                        <p>El Parque Ecoturístico undefined es una opción fantástica para tener una aventura natural en Tabasco. Su calificación promedio es de 4.2 respaldada por más de 33visitantes, así que no tenemos duda de que este lugar pertenece a la lista de los Parque Ecoturístico mejor rankeados de Tabasco y que se trata de uno de los principales atractivos naturales en la región. Así que ya sabes... ¿ganas de naturaleza?... pues el Parque Ecoturístico undefined es una grandísima opción.</p>
                        <h3><b>¿Cómo llegar al Centro Ecoturístico Vivero Santa Fe? </b></h3>
                            <p>El Parque Ecoturístico se ubica enBlvd. Cnel. Enrique Carrola Antuna 909, Ciénega, 34090 Durango, Dgo.. Para ir a este parque, simplemente basta con colocar la dirección en una app de navegación o siguiendo este <a href='https://www.google.com.mx/maps/place/Vivero+Santa+Fe/data=!4m7!3m6!1s0x869bb7e17f5ebdfd:0x329bbfbf57717ca0!8m2!3d24.0199748!4d-104.6558555!16s%2Fg%2F11btmr25_0!19sChIJ_b1ef-G3m4YRoHxxV7-_mzI?authuser=0&hl=es&rclk=1'>Mapa del Parque Ecoturístico Vivero Santa Fe</a></p>
                        <h3><b>¿Cuáles son los contactos del Parque de Ecoturismo Vivero Santa Fe?</b></h3>
                            <p>Los contactos disponibles del Centro de Ecoturismo Vivero Santa Fe son: </p>
                            <ul>
                                <li><b>Teléfono:</b>618 235 9375</li>
                                <li><b>SitioWeb:</b>web no disponible</li>                                
                            </ul>
                        <h3><b>¿En qué horarios y días se puede visitar el Parque Ecoturístico Vivero Santa Fe?</b></h3>
                            <p>Los horarios oficiales del Parque Ecoturístico Vivero Santa Fe son los siguientes:</p>                       
                            <ul>
                                <li>Lunes de 09:30 a 19:30</li>
                                <li>Martes de 09:30 a 19:30</li>
                                <li>Miércoles de 09:30 a 19:30</li>
                                <li>Jueves de 09:30 a 19:30</li>
                                <li>Viernes de 09:30 a 19:30</li>
                                <li>Sábado de 09:30 a 19:30</li>
                                <li>Domingo de 09:30 a 19:30</li>
                            </ul>
                            <p>Aunque cuentes ya con los horarios oficiales de apertura de este lugar, siempre te recomendamos que antes de ir eches un ojito a sus redes sociales y vías de contacto, para asegurarte de que no hayan tenido algún cambio logístico de última hora</p>                 
                    
'@
$ws.Range("T2").Value = $t2
$u2 = @'

                    <img src="undefined" alt="Vivero Santa Fe">                
                
'@
$ws.Range("U2").Value = $u2
$ws.Range("V2").Value = 'Añadir sitio web'

# ---- Row 3 (Vivero Las Magnolias) ----
$ws.Range("N3").Value = 'AF1QipPmQSXM5PoyrFj8IU2GdEEKh9C_3Ku5CJxKUUUM=w408-h306-k-no'
$ws.Range("O3").Value = 'Vivero_Las_Magnolias'
$ws.Range("P3").Value = 'https://rumbonaturaleza.com/uploads/2023/03/Vivero_Las_Magnolias.jpg'
$ws.Range("Q3").Value = 'ren "AF1QipPmQSXM5PoyrFj8IU2GdEEKh9C_3Ku5CJxKUUUM=w408-h306-k-no" "Vivero_Las_Magnolias.jpg"'
$r3 = @'

                    <p> ¿Estás buscando los mejores Parques Ecoturísticos en Tabasco? ¡Estás en el lugar correcto! Pues en este artículo vamos a presentarte cuáles son los  Parque Ecoturístico que han sido mejor evaluados en este estado. 
 Para esto, realizamos consultas en un montón de fuentes oficiales, redes sociales, rankings e incluso entrevistas para poder determinar cuáles son los  Parque Ecoturístico que mejor calificación han recibido en Tabasco durante los últimos años. 
 Con esta prueba social como respaldo, hoy te daremos los Centro de Ecoturismo mejor calificados y te compartiremos su ubicación, medios oficiales de contacto, horarios y cómo llegar hasta ellos, junto con la calificación promedio con la que cuenta cada lugar. 
 Así que prepárate y ¡a disfrutar del ecoturismo en Tabasco!</p>                    
                    
'@
$ws.Range("R3").Value = $r3
$s3 = @'

                    <p> ¿Estás buscando los mejores Parques Ecoturísticos en Tabasco? ¡Estás en el lugar correcto! Pues en este artículo vamos a presentarte cuáles son los  Parque Ecoturístico que han sido mejor evaluados en este estado. 
 Para esto, realizamos consultas en un montón de fuentes oficiales, redes sociales, rankings e incluso entrevistas para poder determinar cuáles son los  Parque Ecoturístico que mejor calificación han recibido en Tabasco durante los últimos años. 
 Con esta prueba social como respaldo, hoy te daremos los Centro de Ecoturismo mejor calificados y te compartiremos su ubicación, medios oficiales de contacto, horarios y cómo llegar hasta ellos, junto con la calificación promedio con la que cuenta cada lugar. 
 Así que prepárate y ¡a disfrutar del ecoturismo en Tabasco!</p>                    
                    
'@
$ws.Range("S3").Value = $s3
$t3 = @'

                    <h2><b>Parque Ecoturístico Vivero Las Magnolias</b></h2>
                        <p>El Parque Ecoturístico undefined es una opción fantástica para tener una aventura natural en Tabasco. Su calificación promedio es de 4.7 respaldada por más de 106visitantes, así que no tenemos duda de que este lugar pertenece a la lista de los Parque Ecoturístico mejor rankeados de Tabasco y que se trata de uno de los principales atractivos naturales en la región. Así que ya sabes... ¿ganas de naturaleza?... pues el Parque Ecoturístico undefined es una grandísima opción.</p>
                        <h3><b>¿Cómo llegar al Parque Natural Vivero Las Magnolias? </b></h3>
                            <p>El Parque Ecoturístico se ubica enSauca 103, Jardines de Durango, 34200 Durango, Dgo.. Para encontrar esta dirección, puedes usar un navegador como el waze o este <a href='https://www.google.com.mx/maps/place/Vivero+Las+Magnolias/data=!4m7!3m6!1s0x869bb7b871872b85:0x65cc8cefdc6b7345!8m2!3d24.0374832!4d-104.6351623!16s%2Fg%2F11b6j5b_lw!19sChIJhSuHcbi3m4YRRXNr3O-MzGU?authuser=0&hl=es&rclk=1'>Mapa del Parque Ecoturístico Vivero Las Magnolias</a></p>
                        <h3><b>¿Cuáles son los contactos del Centro Ecoturístico Vivero Las Magnolias?</b></h3>
                            <p>Los contactos disponibles del Parque Ecológico Vivero Las Magnolias son: </p>
                            <ul>
                                <li><b>Teléfono:</b>618 129 7673</li>
                                <li><b>SitioWeb:</b>web no disponible</li>                                
                            </ul>
                        <h3><b>¿En qué horarios y días se puede visitar el Parque Ecoturístico Vivero Las Magnolias?</b></h3>
                            <p>Los horarios oficiales del Centro Ecoturístico Vivero Las Magnolias son los siguientes:</p>                       
                            <ul>
                                <li>Lunes de 09:00 a 18:00</li>
                                <li>Martes de 09:00 a 18:00</li>
                                <li>Miércoles de 09:00 a 18:00</li>
                                <li>Jueves de 09:00 a 18:00</li>
                                <li>Viernes de 09:00 a 18:00</li>
                                <li>Sábado de 09:00 a 18:00</li>
                                <li>Domingo de 09:00 a 17:00</li>
                            </ul>
                            <p>Aunque estos horarios están oficialmente vigentes, siempre es bueno consultar sus sitios de contacto y redes oficiales antes de visitarlos, por cualquier cambio extraordinario que pudieran tener.</p>                 
                    
'@
$ws.Range("T3").Value = $t3
$u3 = @'

                    <img src="undefined" alt="Vivero Las Magnolias">                
                
'@
$ws.Range("U3").Value = $u3
$ws.Range("V3").Value = 'Añadir sitio web'

